$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12: average of |S*|/n (column J) across the 10 data rows
$ws.Range("J12").Formula = "=AVERAGE(J2:J11)"

# Row 14: Average of SW(S*)/SW(OPT)  -> average of column N
$ws.Range("A14").Value = "Average of SW(S*)/SW(OPT)"
$ws.Range("B14").Formula = "=AVERAGE(N2:N11)"

# Row 15: Average of SC(S*)/SC(OPT) -> average of column Z
$ws.Range("A15").Value = "Average of SC(S*)/SC(OPT)"
$ws.Range("B15").Formula = "=AVERAGE(Z2:Z11)"

# Row 16: Worst of SW(S*)/SW(OPT) -> min of column N
$ws.Range("A16").Value = "Worst of SW(S*)/SW(OPT)"
$ws.Range("B16").Formula = "=MIN(N2:N11)"

# Row 17: Worst of SC(S*)/SC(OPT) -> max of column Z
$ws.Range("A17").Value = "Worst of SC(S*)/SC(OPT)"
$ws.Range("B17").Formula = "=MAX(Z2:Z11)"

# Style B14 fully (bold, size 12, vertical-centered), then propagate the
# exact same formatting to B15:B17 via copy/paste-special so that no stray
# intermediate cell styles are generated.
$b14 = $ws.Range("B14")
$b14.VerticalAlignment = -4108   # xlCenter
$b14.Font.Bold = $true
$b14.Font.Size = 12

$b14.Copy()
$ws.Range("B15:B17").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Row heights to match the taller (12pt) font used in rows 14-17
$ws.Rows.Item(14).RowHeight = 15.6
$ws.Rows.Item(15).RowHeight = 15.6
$ws.Rows.Item(16).RowHeight = 15.6
$ws.Rows.Item(17).RowHeight = 15.6

# Page setup
$ws.PageSetup.PaperSize = 9       # xlPaperA4
$ws.PageSetup.Orientation = 1     # xlPortrait

# Final selection, as left by the author when they saved the file
$ws.Range("A14:B17").Select()
